# Generate Report for Handback
#
# The f61bffeb-d3a3-4c20-ac5f-803e356fc9e4.md file has come back from
# handback processing: its status flips from "Ready for handoff" to
# "Handed back: in sync with en-US" on every sheet, and the zh-cn / de-de
# detail sheets gain their "Latest Target File" / "Latest Handback File"
# values (mirroring the already-handed-off source/file columns) plus a
# real "Latest Handback DateTime" (replacing the 0001-01-01 sentinel).

$wb = $excel.ActiveWorkbook

$HandedBack = "Handed back: in sync with en-US"

# Hyperlink blue used throughout this workbook for "file link" cells
# (style s="1" in the original sheets: underline + RGB 6495ED).
$LinkColor = 15570276  # BGR-encoded OLE color for hex 6495ED

function Set-LinkStyle($range) {
    $range.Font.Underline = 2   # xlUnderlineStyleSingle
    $range.Font.Color = $LinkColor
}

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns for the f61bffeb row
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B5").Value = $HandedBack
$wsOverview.Range("C5").Value = $HandedBack

# ---------------------------------------------------------------------
# zh-cn detail sheet, row 5 (f61bffeb-d3a3-4c20-ac5f-803e356fc9e4.md)
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C5").Value = $HandedBack

$wsZh.Range("F5").Value = "f61bffeb-d3a3-4c20-ac5f-803e356fc9e4.md"
$wsZh.Hyperlinks.Add(
    $wsZh.Range("F5"),
    "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/b82e884638079d8edf522f92ce893692ed463b29/e2e/f61bffeb-d3a3-4c20-ac5f-803e356fc9e4.md",
    "",
    "",
    "f61bffeb-d3a3-4c20-ac5f-803e356fc9e4.md"
) | Out-Null
Set-LinkStyle $wsZh.Range("F5")

$wsZh.Range("G5").Value = "f61bffeb-d3a3-4c20-ac5f-803e356fc9e4.c4d1c593307b6ce628828519568ea3167bfb8d09.zh-cn.xlf"
$wsZh.Hyperlinks.Add(
    $wsZh.Range("G5"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b82e884638079d8edf522f92ce893692ed463b29/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/f61bffeb-d3a3-4c20-ac5f-803e356fc9e4.c4d1c593307b6ce628828519568ea3167bfb8d09.zh-cn.xlf",
    "",
    "",
    "f61bffeb-d3a3-4c20-ac5f-803e356fc9e4.c4d1c593307b6ce628828519568ea3167bfb8d09.zh-cn.xlf"
) | Out-Null
Set-LinkStyle $wsZh.Range("G5")

$wsZh.Range("H5").Value = "2016-03-23 09:46:42"

# ---------------------------------------------------------------------
# de-de detail sheet, row 5 (f61bffeb-d3a3-4c20-ac5f-803e356fc9e4.md)
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C5").Value = $HandedBack

$wsDe.Range("F5").Value = "f61bffeb-d3a3-4c20-ac5f-803e356fc9e4.md"
$wsDe.Hyperlinks.Add(
    $wsDe.Range("F5"),
    "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/4160052fd5bf09f56a254d93a7c47aef5f36fee3/e2e/f61bffeb-d3a3-4c20-ac5f-803e356fc9e4.md",
    "",
    "",
    "f61bffeb-d3a3-4c20-ac5f-803e356fc9e4.md"
) | Out-Null
Set-LinkStyle $wsDe.Range("F5")

$wsDe.Range("G5").Value = "f61bffeb-d3a3-4c20-ac5f-803e356fc9e4.c4d1c593307b6ce628828519568ea3167bfb8d09.de-de.xlf"
$wsDe.Hyperlinks.Add(
    $wsDe.Range("G5"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4160052fd5bf09f56a254d93a7c47aef5f36fee3/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/f61bffeb-d3a3-4c20-ac5f-803e356fc9e4.c4d1c593307b6ce628828519568ea3167bfb8d09.de-de.xlf",
    "",
    "",
    "f61bffeb-d3a3-4c20-ac5f-803e356fc9e4.c4d1c593307b6ce628828519568ea3167bfb8d09.de-de.xlf"
) | Out-Null
Set-LinkStyle $wsDe.Range("G5")

$wsDe.Range("H5").Value = "2016-03-23 09:46:56"
